$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "51.752.64"
Set-TextValue $ws.Range("E2") "  -1.15%  "
Set-TextValue $ws.Range("D3") "2.784.08"
Set-TextValue $ws.Range("E3") "  -1.68%  "
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("D5") "357.29"
Set-TextValue $ws.Range("E5") "  +0.27%  "
Set-TextValue $ws.Range("D6") "109.23"
Set-TextValue $ws.Range("E6") "  -2.90%  "
Set-TextValue $ws.Range("D7") "0.555"
Set-TextValue $ws.Range("E7") "  -3.21%  "
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("D9") "0.588"
Set-TextValue $ws.Range("E9") "  -2.27%  "
Set-TextValue $ws.Range("D10") "39.72"
Set-TextValue $ws.Range("E10") "  -3.42%  "
Set-TextValue $ws.Range("E11") "  +3.67%  "
Set-TextValue $ws.Range("D12") "0.0847"
Set-TextValue $ws.Range("E12") "  -2.12%  "
Set-TextValue $ws.Range("D13") "19.57"
Set-TextValue $ws.Range("E13") "  -1.93%  "
Set-TextValue $ws.Range("D14") "7.60"
Set-TextValue $ws.Range("E14") "  -2.38%  "
Set-TextValue $ws.Range("D15") "3.218.47"
Set-TextValue $ws.Range("E15") "  -1.56%  "
Set-TextValue $ws.Range("D16") "2.777.35"
Set-TextValue $ws.Range("E16") "  -1.83%  "
Set-TextValue $ws.Range("D17") "0.934"
Set-TextValue $ws.Range("E17") "  +0.44%  "
Set-TextValue $ws.Range("D18") "51.701.47"
Set-TextValue $ws.Range("E18") "  -0.76%  "
Set-TextValue $ws.Range("D19") "7.56"
Set-TextValue $ws.Range("E19") "  +0.60%  "
Set-TextValue $ws.Range("D20") "3.08"
Set-TextValue $ws.Range("E20") "  -3.64%  "
Set-TextValue $ws.Range("D21") "13.23"
Set-TextValue $ws.Range("E21") "  -2.15%  "
Set-TextValue $ws.Range("D22") "0.0₃0970"
Set-TextValue $ws.Range("E22") "  -2.75%  "
Set-TextValue $ws.Range("D23") "70.26"
Set-TextValue $ws.Range("E23") "  -0.55%  "
Set-TextValue $ws.Range("D24") "267.85"
Set-TextValue $ws.Range("E24") "  -1.53%  "
Set-TextValue $ws.Range("E25") "  -2.45%  "
Set-TextValue $ws.Range("D26") "26.38"
Set-TextValue $ws.Range("E26") "  -2.35%  "
Set-TextValue $ws.Range("B27") "Kaspa"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D27") "0.168"
Set-TextValue $ws.Range("E27") "  +16.50%  "
Set-TextValue $ws.Range("B28") "Dai"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D28") "1.00"
Set-TextValue $ws.Range("E28") "  -0.04%  "
Set-TextValue $ws.Range("D29") "10.21"
Set-TextValue $ws.Range("E29") "  -1.43%  "
Set-TextValue $ws.Range("E30") "  -3.88%  "
Set-TextValue $ws.Range("D31") "6.16"
Set-TextValue $ws.Range("E31") "  +3.54%  "
Set-TextValue $ws.Range("B32") "InjectiveProtocol"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D32") "35.05"
Set-TextValue $ws.Range("E32") "  -0.20%  "
Set-TextValue $ws.Range("B33") "OKB"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D33") "51.93"
Set-TextValue $ws.Range("E33") "  -1.36%  "
Set-TextValue $ws.Range("E34") "  -8.83%  "
Set-TextValue $ws.Range("D35") "0.0836"
Set-TextValue $ws.Range("E35") "  -2.59%  "
Set-TextValue $ws.Range("D36") "5.20"
Set-TextValue $ws.Range("E36") "  -7.20%  "
Set-TextValue $ws.Range("D37") "0.999"
Set-TextValue $ws.Range("E37") "  -0.02%  "
Set-TextValue $ws.Range("D38") "18.88"
Set-TextValue $ws.Range("E38") "  +2.38%  "
Set-TextValue $ws.Range("E39") "  -4.20%  "
Set-TextValue $ws.Range("E40") "  -4.25%  "
Set-TextValue $ws.Range("E41") "  +0.38%  "
Set-TextValue $ws.Range("E42") "  -2.99%  "
Set-TextValue $ws.Range("E43") "  -3.02%  "
Set-TextValue $ws.Range("D44") "119.83"
Set-TextValue $ws.Range("E44") "  -5.95%  "
Set-TextValue $ws.Range("D45") "21.75"
Set-TextValue $ws.Range("E45") "  -6.68%  "
Set-TextValue $ws.Range("D46") "2.086.90"
Set-TextValue $ws.Range("E46") "  -0.11%  "
Set-TextValue $ws.Range("D47") "3.27"
Set-TextValue $ws.Range("E47") "  -2.86%  "
Set-TextValue $ws.Range("E48") "  -0.08%  "
Set-TextValue $ws.Range("D49") "0.948"
Set-TextValue $ws.Range("E49") "  -2.62%  "
Set-TextValue $ws.Range("D50") "5.55"
Set-TextValue $ws.Range("E50") "  -6.70%  "
Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.191"
Set-TextValue $ws.Range("E51") "  -1.87%  "
